# Deploy the implementation guide: refresh the generated ValueSet metadata
# sheet with the new publish status/date from the latest IG build.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date: 2023-05-12T12:33:13+00:00 -> 2023-08-01T16:12:28+00:00
$ws.Range("B8").Value = "2023-08-01T16:12:28+00:00"
